$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target value still parses as a plain number need the
# column kept as Text so Excel does not silently reformat/round them
# (e.g. "7.40" -> 7.4, "1.00" -> 1). Force NumberFormat = "@" on just
# those cells before writing the new text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "624.06"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.65"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.40"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.62"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.13"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.19"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "462.70"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.647"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.07"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000133"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.22"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.171"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.54"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.95"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.42"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.41"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0926"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "179.53"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.996"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.65"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.18"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.907"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.95"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.71"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.86"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.267"

# Remaining changed cells (already non-numeric text, or intentionally
# numeric-looking-but-unparseable strings like "69.075.89") can be set
# directly.
$ws.Range("D2").Value = "69.075.89"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").Value = "3.584.35"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("E6").Value = "  +6.39%  "
$ws.Range("D7").Value = "3.589.23"
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("E10").Value = "  +7.96%  "
$ws.Range("E11").Value = "  +6.88%  "
$ws.Range("E12").Value = "  +5.06%  "
$ws.Range("E13").Value = "  +4.92%  "
$ws.Range("E14").Value = "  +7.21%  "
$ws.Range("D15").Value = "4.194.95"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "69.556.16"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("D17").Value = "3.576.01"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +6.65%  "
$ws.Range("E20").Value = "  +7.72%  "
$ws.Range("E21").Value = "  +13.04%  "
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("E23").Value = "  +4.15%  "
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("E25").Value = "  +5.55%  "
$ws.Range("E26").Value = "  +5.44%  "
$ws.Range("D27").Value = "3.724.32"
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +11.49%  "
$ws.Range("E30").Value = "  +9.04%  "
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("E32").Value = "  +5.50%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E33").Value = "  +7.15%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E35").Value = "  +5.57%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("D37").Value = "3.581.23"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("E38").Value = "  +5.36%  "
$ws.Range("E39").Value = "  +9.54%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E41").Value = "  +6.43%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E42").Value = "  +5.37%  "
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("E45").Value = "  +17.33%  "
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("E47").Value = "  +9.86%  "
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  +6.93%  "
$ws.Range("E50").Value = "  +4.41%  "
$ws.Range("E51").Value = "  +9.88%  "
